$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (so existing B..I shift to C..J)
$ws.Columns("B").Insert()

# Set header for new column B
$ws.Cells.Item(1, 2).Value = "status_label"

# Set "rouge" value for each data row (rows 2-6) in column B
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 2).Value = "rouge"
}
